$wb = $excel.ActiveWorkbook

# Rename the "SwateTemplateMetadata" sheet to "isa_template"
$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metaSheet.Name = "isa_template"

# Update the active selection on that sheet to C9 (from B4)
$metaSheet.Activate()
$metaSheet.Range("C9").Select()
